$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.532.35'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '2.338.78'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  +2.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.97'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.65%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.457'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0972'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = '2.688.20'
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.858'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").Value = '2.343.96'
$ws.Range("E18").Value = '  +2.90%  '
$ws.Range("D19").Value = '43.552.75'
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +13.54%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("E31").Value = '  +5.84%  '
$ws.Range("E32").Value = '  -6.94%  '
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("E34").Value = '  +4.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0690'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("E37").Value = '  +9.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0254'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.75%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.72%  '
$ws.Range("E44").Value = '  +9.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0949'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.448.51'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("B50").Value = 'Celestia'
$ws.Range("C50").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.58%  '
$ws.Range("B51").Value = 'TerraClassic'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000206'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -14.04%  '
